$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume/change (E) columns
# Cells whose new numeric-looking text would otherwise be auto-converted
# to a Number by Excel (dropping trailing zeros) are forced to Text first.

$ws.Range("D2").Value = '63.454.02'
$ws.Range("E2").Value = '  +0.96%  '

$ws.Range("D3").Value = '2.679.30'
$ws.Range("E3").Value = '  +4.34%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '614.70'
$ws.Range("E5").Value = '  +5.57%  '

$ws.Range("D6").Value = '144.08'
$ws.Range("E6").Value = '  +0.33%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("E8").Value = '  -0.10%  '

$ws.Range("D9").Value = '2.677.27'
$ws.Range("E9").Value = '  +4.29%  '

$ws.Range("E10").Value = '  +1.22%  '

$ws.Range("D11").Value = '5.63'
$ws.Range("E11").Value = '  +1.08%  '

$ws.Range("E12").Value = '  +0.69%  '

$ws.Range("D13").Value = '0.363'
$ws.Range("E13").Value = '  +4.47%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.40'
$ws.Range("E14").Value = '  +1.73%  '

$ws.Range("D15").Value = '3.158.26'
$ws.Range("E15").Value = '  +4.18%  '

$ws.Range("D16").Value = '63.312.09'
$ws.Range("E16").Value = '  +0.91%  '

$ws.Range("D17").Value = '0.0000145'
$ws.Range("E17").Value = '  +0.94%  '

$ws.Range("D18").Value = '2.680.46'
$ws.Range("E18").Value = '  +3.65%  '

$ws.Range("E19").Value = '  +3.94%  '

$ws.Range("D20").Value = '343.43'
$ws.Range("E20").Value = '  +1.03%  '

$ws.Range("E21").Value = '  +2.30%  '

$ws.Range("E22").Value = '  +4.13%  '

$ws.Range("E23").Value = '  -0.07%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '67.40'
$ws.Range("E24").Value = '  -0.36%  '

$ws.Range("E25").Value = '  +4.44%  '

$ws.Range("E26").Value = '  -2.15%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.70'
$ws.Range("E27").Value = '  +5.90%  '

$ws.Range("E28").Value = '  +0.04%  '

$ws.Range("D29").Value = '543.03'
$ws.Range("E29").Value = '  +18.38%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.15%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.90'
$ws.Range("E31").Value = '  -0.77%  '

$ws.Range("E32").Value = '  +7.54%  '

$ws.Range("E33").Value = '  +8.95%  '

$ws.Range("E34").Value = '  +2.17%  '

$ws.Range("D35").Value = '172.55'
$ws.Range("E35").Value = '  -2.31%  '

$ws.Range("E36").Value = '  +14.76%  '

$ws.Range("E37").Value = '  +2.02%  '

$ws.Range("E38").Value = '  -0.09%  '

$ws.Range("D39").Value = '19.27'
$ws.Range("E39").Value = '  +2.54%  '

$ws.Range("E40").Value = '  +10.74%  '

$ws.Range("D41").Value = '176.44'
$ws.Range("E41").Value = '  +12.24%  '

$ws.Range("E42").Value = '  -0.01%  '

$ws.Range("E43").Value = '  +2.40%  '

$ws.Range("D44").Value = '22.31'
$ws.Range("E44").Value = '  +5.55%  '

$ws.Range("E45").Value = '  +7.46%  '

$ws.Range("E46").Value = '  +1.02%  '

$ws.Range("E47").Value = '  +3.03%  '

$ws.Range("D48").Value = '0.0966'
$ws.Range("E48").Value = '  +0.89%  '

$ws.Range("D49").Value = '18.96'
$ws.Range("E49").Value = '  +5.43%  '

$ws.Range("E50").Value = '  +5.71%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '11.30'
$ws.Range("E51").Value = '  -0.83%  '
